$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19: Currò Gaetano / Spilii / Scatole / Mt. 1 ---
$ws.Range("A19").Value = 43228
$ws.Range("B19").Value = "Currò Gaetano"
$ws.Range("C19").Value = "Spilii"
$ws.Range("D19").Value = "Scatole"
$ws.Range("E19").Value = 1

# --- Row 20: Bertolotti Daniela / Adesivo Leggero nero / Mt. 80 ---
$ws.Range("A20").Value = 43231
$ws.Range("B20").Value = "Bertolotti Daniela"
$ws.Range("C20").Value = "Adesivo Leggero nero"
$ws.Range("D20").Value = "Mt."
$ws.Range("E20").Value = 80

# Apply the same cell formatting used by the rest of the table (row 18 is
# the last pre-existing data row) to the two newly appended rows.
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A18:E18").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
